# Add a new worksheet "with separators" right after the existing "data" sheet,
# populate it with sample data that uses " | " and " ! " style separators plus
# a lone-space value, then leave the selection positioned below the data —
# mirroring the commit that adds regression-test data for the
# table-properties-transformer fix.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$newSheet = $wb.Worksheets.Add([Type]::Missing, $dataSheet)
$newSheet.Name = "with separators"

$newSheet.Range("A1").Value = "A | B"
$newSheet.Range("B1").Value = "C ! D"
$newSheet.Range("A2").Value = " "
$newSheet.Range("B2").Value = " "

$null = $newSheet.Range("B3").Select()
